$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VLRS")

$ws.Range("D8").Value = 268000
$ws.Range("E8").Value = 263400
$ws.Range("F8").Value = 270000
$ws.Range("G8").Value = 212000
$ws.Range("H8").Value = 194900
$ws.Range("I8").Value = 204700
$ws.Range("J8").Value = 206600
$ws.Range("D17").Value = 224800
$ws.Range("E17").Value = 225300
$ws.Range("F17").Value = 214100
$ws.Range("G17").Value = 194000
$ws.Range("H17").Value = 177000
$ws.Range("I17").Value = 182700
$ws.Range("J17").Value = 187900
$ws.Range("D18").Value = 43200
$ws.Range("E18").Value = 38000
$ws.Range("F18").Value = 55800
$ws.Range("G18").Value = 18000
$ws.Range("H18").Value = 17900
$ws.Range("I18").Value = 22000
$ws.Range("J18").Value = 18700
$ws.Range("E20").Value = 9400
$ws.Range("F20").Value = 29300
$ws.Range("G20").Value = 7900
$ws.Range("H20").Value = 4700
$ws.Range("I20").Value = 17300
$ws.Range("J20").Value = 5900
$ws.Range("D21").Value = 51300
$ws.Range("E21").Value = 53000
$ws.Range("F21").Value = 91400
$ws.Range("G21").Value = 32400
$ws.Range("H21").Value = 28000
$ws.Range("I21").Value = 46400
$ws.Range("J21").Value = 29000
$ws.Range("D23").Value = 45100
$ws.Range("E23").Value = 47400
$ws.Range("F23").Value = 85100
$ws.Range("G23").Value = 26000
$ws.Range("H23").Value = 22600
$ws.Range("I23").Value = 39300
$ws.Range("J23").Value = 24500
$ws.Range("D24").Value = 14000
$ws.Range("E24").Value = 13600
$ws.Range("F24").Value = 25600
$ws.Range("G24").Value = 7800
$ws.Range("H24").Value = 6800
$ws.Range("J24").Value = 6600
$ws.Range("D26").Value = 31100
$ws.Range("E26").Value = 33800
$ws.Range("F26").Value = 59600
$ws.Range("G26").Value = 18200
$ws.Range("H26").Value = 15800
$ws.Range("I26").Value = 36400
$ws.Range("J26").Value = 18000
$ws.Range("D27").Value = 31100
$ws.Range("E27").Value = 33800
$ws.Range("F27").Value = 59600
$ws.Range("G27").Value = 18200
$ws.Range("H27").Value = 15800
$ws.Range("I27").Value = 36400
$ws.Range("J27").Value = 18000
$ws.Range("E32").Value = -9400
$ws.Range("F32").Value = -29300
$ws.Range("G32").Value = -7900
$ws.Range("H32").Value = -4700
$ws.Range("I32").Value = -17300
$ws.Range("J32").Value = -5900
$ws.Range("D33").Value = 31100
$ws.Range("E33").Value = 33800
$ws.Range("F33").Value = 59600
$ws.Range("G33").Value = 18200
$ws.Range("H33").Value = 15800
$ws.Range("I33").Value = 36400
$ws.Range("J33").Value = 18000
$ws.Range("D35").Value = 31100
$ws.Range("E35").Value = 33800
$ws.Range("F35").Value = 59600
$ws.Range("G35").Value = 18200
$ws.Range("H35").Value = 15800
$ws.Range("I35").Value = 36400
$ws.Range("J35").Value = 18000
$ws.Range("D41").Value = 329300
$ws.Range("E41").Value = 266700
$ws.Range("F41").Value = 228000
$ws.Range("G41").Value = 208300
$ws.Range("H41").Value = 163200
$ws.Range("I41").Value = 160700
$ws.Range("J41").Value = 93800
$ws.Range("I42").Value = 31400
$ws.Range("D43").Value = 24300
$ws.Range("E43").Value = 24000
$ws.Range("F43").Value = 15300
$ws.Range("G43").Value = 14800
$ws.Range("H43").Value = 14800
$ws.Range("I43").Value = 35300
$ws.Range("J43").Value = 28300
$ws.Range("D44").Value = 8800
$ws.Range("E44").Value = 8400
$ws.Range("F44").Value = 8100
$ws.Range("G44").Value = 8200
$ws.Range("H44").Value = 7100
$ws.Range("I44").Value = 14400
$ws.Range("J44").Value = 6600
$ws.Range("D45").Value = 88500
$ws.Range("E45").Value = 75300
$ws.Range("F45").Value = 54900
$ws.Range("G45").Value = 52100
$ws.Range("H45").Value = 45800
$ws.Range("I45").Value = 116800
$ws.Range("J45").Value = 45300
$ws.Range("D46").Value = 450900
$ws.Range("E46").Value = 374500
$ws.Range("F46").Value = 306200
$ws.Range("G46").Value = 283400
$ws.Range("H46").Value = 230900
$ws.Range("I46").Value = 190800
$ws.Range("J46").Value = 174100
$ws.Range("I47").Value = 183400
$ws.Range("D48").Value = 113100
$ws.Range("E48").Value = 131900
$ws.Range("F48").Value = 117600
$ws.Range("G48").Value = 124700
$ws.Range("H48").Value = 113300
$ws.Range("I48").Value = 230000
$ws.Range("J48").Value = 103000
$ws.Range("D49").Value = 5000
$ws.Range("E49").Value = 4900
$ws.Range("F49").Value = 4000
$ws.Range("G49").Value = 3500
$ws.Range("H49").Value = 3600
$ws.Range("I49").Value = 7500
$ws.Range("J49").Value = 3100
$ws.Range("D52").Value = 289700
$ws.Range("E52").Value = 278000
$ws.Range("F52").Value = 283000
$ws.Range("G52").Value = 247600
$ws.Range("H52").Value = 225500
$ws.Range("I52").Value = 205200
$ws.Range("J52").Value = 169800
$ws.Range("D54").Value = 858800
$ws.Range("E54").Value = 789300
$ws.Range("F54").Value = 710800
$ws.Range("G54").Value = 659200
$ws.Range("H54").Value = 573300
$ws.Range("I54").Value = 512300
$ws.Range("J54").Value = 450000
$ws.Range("D57").Value = 245400
$ws.Range("E57").Value = 40400
$ws.Range("F57").Value = 35900
$ws.Range("G57").Value = 27100
$ws.Range("H57").Value = 13700
$ws.Range("I57").Value = 26200
$ws.Range("J57").Value = 27500
$ws.Range("D58").Value = 46000
$ws.Range("E58").Value = 70500
$ws.Range("F58").Value = 58900
$ws.Range("G58").Value = 65900
$ws.Range("H58").Value = 43800
$ws.Range("I58").Value = 42300
$ws.Range("J58").Value = 13800
$ws.Range("D59").Value = 108900
$ws.Range("E59").Value = 256400
$ws.Range("F59").Value = 256400
$ws.Range("G59").Value = 273800
$ws.Range("H59").Value = 234400
$ws.Range("I59").Value = 271600
$ws.Range("J59").Value = 166700
$ws.Range("D60").Value = 400300
$ws.Range("E60").Value = 367400
$ws.Range("F60").Value = 351100
$ws.Range("G60").Value = 366800
$ws.Range("H60").Value = 291800
$ws.Range("I60").Value = 246600
$ws.Range("J60").Value = 208000
$ws.Range("E61").Value = 11400
$ws.Range("F61").Value = 16000
$ws.Range("G61").Value = 9000
$ws.Range("H61").Value = 20900
$ws.Range("I61").Value = 22000
$ws.Range("J61").Value = 32100
$ws.Range("D62").Value = 71400
$ws.Range("E62").Value = 57600
$ws.Range("F62").Value = 23200
$ws.Range("G62").Value = 13700
$ws.Range("H62").Value = 12000
$ws.Range("I62").Value = 17600
$ws.Range("J62").Value = 10900
$ws.Range("D66").Value = 471700
$ws.Range("E66").Value = 436300
$ws.Range("F66").Value = 390300
$ws.Range("G66").Value = 389500
$ws.Range("H66").Value = 324700
$ws.Range("I66").Value = 281100
$ws.Range("J66").Value = 251000
$ws.Range("D72").Value = 155700
$ws.Range("E72").Value = 126500
$ws.Range("F72").Value = 92700
$ws.Range("G72").Value = 33100
$ws.Range("H72").Value = 14900
$ws.Range("I72").Value = -3800
$ws.Range("J72").Value = -37300
$ws.Range("D76").Value = 387100
$ws.Range("E76").Value = 353000
$ws.Range("F76").Value = 320500
$ws.Range("G76").Value = 269700
$ws.Range("H76").Value = 248600
$ws.Range("I76").Value = 231200
$ws.Range("J76").Value = 199000
$ws.Range("D81").Value = 31100
$ws.Range("E81").Value = 33800
$ws.Range("F81").Value = 59600
$ws.Range("G81").Value = 18200
$ws.Range("H81").Value = 15800
$ws.Range("I81").Value = 36400
$ws.Range("J81").Value = 18000
$ws.Range("D83").Value = 6200
$ws.Range("E83").Value = 5600
$ws.Range("F83").Value = 6200
$ws.Range("G83").Value = 6500
$ws.Range("H83").Value = 5300
$ws.Range("I83").Value = 7100
$ws.Range("J83").Value = 4500
$ws.Range("D89").Value = 68700
$ws.Range("E89").Value = 48100
$ws.Range("F89").Value = 12600
$ws.Range("G89").Value = 49000
$ws.Range("H89").Value = 49100
$ws.Range("I89").Value = 24300
$ws.Range("D91").Value = -10200
$ws.Range("E91").Value = -26900
$ws.Range("F91").Value = -12600
$ws.Range("G91").Value = -22600
$ws.Range("H91").Value = -10500
$ws.Range("I91").Value = -25400
$ws.Range("J91").Value = -19000
$ws.Range("D94").Value = 22500
$ws.Range("E94").Value = -18400
$ws.Range("F94").Value = 4500
$ws.Range("G94").Value = -14500
$ws.Range("H94").Value = -2600
$ws.Range("I94").Value = -19300
$ws.Range("J94").Value = -19100
$ws.Range("D100").Value = -28400
$ws.Range("E100").Value = 6500
$ws.Range("F100").Value = -9100
$ws.Range("G100").Value = 7800
$ws.Range("H100").Value = -1900
$ws.Range("I100").Value = 12700
$ws.Range("J100").Value = 5000
$ws.Range("F101").Value = 11700
$ws.Range("G101").Value = 2800
$ws.Range("I101").Value = 5600
$ws.Range("D102").Value = 62500
$ws.Range("E102").Value = 38800
$ws.Range("F102").Value = 19600
$ws.Range("G102").Value = 45100
$ws.Range("H102").Value = 46100
$ws.Range("I102").Value = 23300
$ws.Range("J102").Value = -14200
